# Generate Report for Handoff
# Adds a new tracked file ("f843c243-6809-434c-80ed-ff15e90cb237") as row 9
# to the Overview sheet and to each per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$fileId     = "f843c243-6809-434c-80ed-ff15e90cb237"
$mdName     = "$fileId.md"
$commitHash = "4c3edea4343f13c7713e0d09b8069b7f97e3d4d6"
$srcCommit  = "1f6c275ccfcd8e71a1d639e6b1b1c4d0c6fa9b21"

# Hyperlink-blue color used throughout the workbook (RGB 6495ED), expressed
# as a COLORREF (BGR order) for the Font.Color COM property.
$linkColor = 15570276

function Style-AsLink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $linkColor
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A9").Value = $mdName
$ws1.Range("B9").Value = "Ready for handoff"
$ws1.Range("C9").Value = "Ready for handoff"
$ws1.Range("D9").Value = "2016-30-12 12:30:56"

$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdName", "", "", $mdName) | Out-Null
Style-AsLink $ws1.Range("A9")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$zhXlf = "$fileId.$commitHash.zh-cn.xlf"

$ws2.Range("A9").Value = $mdName
$ws2.Range("B9").Value = ".md"
$ws2.Range("C9").Value = "Ready for handoff"
$ws2.Range("D9").Value = $zhXlf
$ws2.Range("E9").Value = "2016-03-12 12:30:53"
$ws2.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H9").Value = "0001-01-01 00:00:00"
$ws2.Range("I9").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdName", "", "", $mdName) | Out-Null
Style-AsLink $ws2.Range("A9")

$ws2.Hyperlinks.Add($ws2.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdName", "", "", ".md") | Out-Null
Style-AsLink $ws2.Range("B9")

$ws2.Hyperlinks.Add($ws2.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf) | Out-Null
Style-AsLink $ws2.Range("D9")

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$deXlf = "$fileId.$commitHash.de-de.xlf"

$ws3.Range("A9").Value = $mdName
$ws3.Range("B9").Value = ".md"
$ws3.Range("C9").Value = "Ready for handoff"
$ws3.Range("D9").Value = $deXlf
$ws3.Range("E9").Value = "2016-03-12 12:30:56"
$ws3.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H9").Value = "0001-01-01 00:00:00"
$ws3.Range("I9").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdName", "", "", $mdName) | Out-Null
Style-AsLink $ws3.Range("A9")

$ws3.Hyperlinks.Add($ws3.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdName", "", "", ".md") | Out-Null
Style-AsLink $ws3.Range("B9")

$ws3.Hyperlinks.Add($ws3.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf) | Out-Null
Style-AsLink $ws3.Range("D9")
